$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '73.872.84'
Set-TextCell 2 5 '  +7.47%  '

Set-TextCell 3 4 '2.628.01'
Set-TextCell 3 5 '  +7.84%  '

Set-TextCell 4 4 '1.00'
Set-TextCell 4 5 '  -0.05%  '

Set-TextCell 5 4 '185.13'
Set-TextCell 5 5 '  +14.60%  '

Set-TextCell 6 4 '583.21'
Set-TextCell 6 5 '  +4.41%  '

Set-TextCell 7 5 '  -0.13%  '

Set-TextCell 8 5 '  +4.38%  '

Set-TextCell 9 4 '0.202'
Set-TextCell 9 5 '  +19.79%  '

Set-TextCell 10 4 '2.627.46'
Set-TextCell 10 5 '  +7.89%  '

Set-TextCell 11 5 '  +0.34%  '

Set-TextCell 12 5 '  +8.25%  '

Set-TextCell 13 5 '  +4.02%  '

Set-TextCell 14 5 '  +6.98%  '

Set-TextCell 15 2 'WrappedliquidstakedEther2.0'
Set-TextCell 15 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 15 4 '3.108.58'
Set-TextCell 15 5 '  +7.74%  '

Set-TextCell 16 2 'WrappedBTC'
Set-TextCell 16 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 16 4 '73.777.56'
Set-TextCell 16 5 '  +7.50%  '

Set-TextCell 17 4 '26.20'
Set-TextCell 17 5 '  +12.94%  '

Set-TextCell 18 4 '2.628.00'
Set-TextCell 18 5 '  +7.82%  '

Set-TextCell 19 4 '9.07'
Set-TextCell 19 5 '  +30.81%  '

Set-TextCell 20 4 '11.86'
Set-TextCell 20 5 '  +12.33%  '

Set-TextCell 21 4 '371.83'
Set-TextCell 21 5 '  +9.58%  '

Set-TextCell 22 5 '  +19.32%  '

Set-TextCell 23 5 '  +6.97%  '

Set-TextCell 24 5 '  +0.15%  '

Set-TextCell 25 4 '69.87'
Set-TextCell 25 5 '  +4.50%  '

Set-TextCell 26 5 '  +11.84%  '

Set-TextCell 27 4 '9.41'
Set-TextCell 27 5 '  +14.79%  '

Set-TextCell 28 4 '2.763.19'
Set-TextCell 28 5 '  +7.76%  '

Set-TextCell 29 4 '1.00'
Set-TextCell 29 5 '  -0.01%  '

Set-TextCell 30 4 '0.0₃0944'
Set-TextCell 30 5 '  +15.18%  '

Set-TextCell 31 4 '523.20'
Set-TextCell 31 5 '  +22.37%  '

Set-TextCell 32 5 '  +20.56%  '

Set-TextCell 33 4 '7.67'
Set-TextCell 33 5 '  +7.67%  '

Set-TextCell 34 5 '  +9.36%  '

Set-TextCell 35 5 '  -0.09%  '

Set-TextCell 36 5 '  +13.76%  '

Set-TextCell 37 4 '160.92'
Set-TextCell 37 5 '  +0.93%  '

Set-TextCell 38 4 '19.18'
Set-TextCell 38 5 '  +6.64%  '

Set-TextCell 39 5 '  +1.49%  '

Set-TextCell 40 5 '  -0.03%  '

Set-TextCell 41 4 '4.92'
Set-TextCell 41 5 '  +13.26%  '

Set-TextCell 42 5 '  +9.77%  '

Set-TextCell 43 5 '  +10.97%  '

Set-TextCell 44 5 '  +23.39%  '

Set-TextCell 45 4 '1.19'
Set-TextCell 45 5 '  +11.11%  '

Set-TextCell 46 5 '  +15.66%  '

Set-TextCell 47 4 '38.99'
Set-TextCell 47 5 '  +4.11%  '

Set-TextCell 48 4 '0.0850'
Set-TextCell 48 5 '  +18.31%  '

Set-TextCell 49 5 '  +9.07%  '

Set-TextCell 50 4 '0.529'
Set-TextCell 50 5 '  +10.16%  '

Set-TextCell 51 4 '20.73'
Set-TextCell 51 5 '  +22.86%  '
